$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 9881
$ws.Range("E3").Value = 2343
$ws.Range("E4").Value = 10950
$ws.Range("E5").Value = 2697
$ws.Range("E6").Value = 19513
$ws.Range("E7").Value = 3655
$ws.Range("E8").Value = 1128
$ws.Range("E9").Value = 5248
$ws.Range("E10").Value = 14597
$ws.Range("E11").Value = 13065
$ws.Range("E12").Value = 16663
$ws.Range("E13").Value = 11807
